$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 26 de Junio de 2020 a las 13:17"

# Apply updated COVID statistics (values refreshed + re-sorted rows)
$ws.Cells.Item(7, 2).Value = 491992
$ws.Cells.Item(7, 3).Value = 822
$ws.Cells.Item(7, 4).Value = 286019
$ws.Cells.Item(7, 5).Value = 190654
$ws.Cells.Item(13, 2).Value = 217724
$ws.Cells.Item(13, 3).Value = 2628
$ws.Cells.Item(13, 4).Value = 177852
$ws.Cells.Item(13, 5).Value = 29633
$ws.Cells.Item(13, 7).Value = 109
$ws.Cells.Item(13, 8).Value = 10239
$ws.Cells.Item(16, 2).Value = 193807
$ws.Cells.Item(16, 3).Value = 22
$ws.Cells.Item(16, 4).Value = 177100
$ws.Cells.Item(16, 5).Value = 7695
$ws.Cells.Item(23, 2).Value = 92784
$ws.Cells.Item(23, 3).Value = 946
$ws.Cells.Item(23, 4).Value = 76072
$ws.Cells.Item(23, 5).Value = 16603
$ws.Cells.Item(23, 7).Value = 3
$ws.Cells.Item(23, 8).Value = 109
$ws.Cells.Item(28, 2).Value = 61106
$ws.Cells.Item(28, 3).Value = 99
$ws.Cells.Item(28, 4).Value = 16918
$ws.Cells.Item(28, 5).Value = 34457
$ws.Cells.Item(28, 7).Value = 5
$ws.Cells.Item(28, 8).Value = 9731
$ws.Cells.Item(32, 2).Value = 51427
$ws.Cells.Item(32, 3).Value = 1240
$ws.Cells.Item(32, 4).Value = 21333
$ws.Cells.Item(32, 5).Value = 27411
$ws.Cells.Item(32, 7).Value = 63
$ws.Cells.Item(32, 8).Value = 2683
$ws.Cells.Item(40, 2).Value = 36034
$ws.Cells.Item(40, 3).Value = 1132
$ws.Cells.Item(40, 4).Value = 19482
$ws.Cells.Item(40, 5).Value = 16399
$ws.Cells.Item(40, 7).Value = 9
$ws.Cells.Item(40, 8).Value = 153
$ws.Cells.Item(41, 1).Value = "Filipinas"
$ws.Cells.Item(41, 2).Value = 34073
$ws.Cells.Item(41, 3).Value = 1004
$ws.Cells.Item(41, 4).Value = 9182
$ws.Cells.Item(41, 5).Value = 23667
$ws.Cells.Item(41, 7).Value = 12
$ws.Cells.Item(41, 8).Value = 1224
$ws.Cells.Item(42, 1).Value = "Polonia"
$ws.Cells.Item(42, 2).Value = 33395
$ws.Cells.Item(42, 3).Value = 276
$ws.Cells.Item(42, 4).Value = 19218
$ws.Cells.Item(42, 5).Value = 12748
$ws.Cells.Item(42, 7).Value = 17
$ws.Cells.Item(42, 8).Value = 1429
$ws.Cells.Item(43, 2).Value = 31486
$ws.Cells.Item(43, 3).Value = 58
$ws.Cells.Item(43, 5).Value = 524
$ws.Cells.Item(43, 7).Value = 4
$ws.Cells.Item(43, 8).Value = 1962
$ws.Cells.Item(48, 1).Value = "Rumania"
$ws.Cells.Item(48, 2).Value = 25697
$ws.Cells.Item(48, 3).Value = 411
$ws.Cells.Item(48, 4).Value = 18181
$ws.Cells.Item(48, 5).Value = 5937
$ws.Cells.Item(48, 7).Value = 14
$ws.Cells.Item(48, 8).Value = 1579
$ws.Cells.Item(49, 1).Value = "Irlanda"
$ws.Cells.Item(49, 2).Value = 25405
$ws.Cells.Item(49, 4).Value = 23364
$ws.Cells.Item(49, 5).Value = 314
$ws.Cells.Item(49, 8).Value = 1727
$ws.Cells.Item(52, 1).Value = "Israel"
$ws.Cells.Item(52, 2).Value = 22638
$ws.Cells.Item(52, 3).Value = 238
$ws.Cells.Item(52, 4).Value = 16589
$ws.Cells.Item(52, 5).Value = 5736
$ws.Cells.Item(52, 7).Value = 4
$ws.Cells.Item(52, 8).Value = 313
$ws.Cells.Item(53, 1).Value = "Nigeria"
$ws.Cells.Item(53, 2).Value = 22614
$ws.Cells.Item(53, 4).Value = 7822
$ws.Cells.Item(53, 5).Value = 14243
$ws.Cells.Item(53, 8).Value = 549
$ws.Cells.Item(56, 2).Value = 17522
$ws.Cells.Item(56, 3).Value = 45
$ws.Cells.Item(56, 4).Value = 16348
$ws.Cells.Item(56, 5).Value = 476
$ws.Cells.Item(67, 1).Value = "Nepal"
$ws.Cells.Item(67, 2).Value = 11755
$ws.Cells.Item(67, 3).Value = 593
$ws.Cells.Item(67, 4).Value = 2698
$ws.Cells.Item(67, 5).Value = 9030
$ws.Cells.Item(67, 7).Value = 1
$ws.Cells.Item(67, 8).Value = 27
$ws.Cells.Item(68, 1).Value = "Marruecos"
$ws.Cells.Item(68, 2).Value = 11465
$ws.Cells.Item(68, 3).Value = 127
$ws.Cells.Item(68, 4).Value = 8560
$ws.Cells.Item(68, 5).Value = 2688
$ws.Cells.Item(68, 8).Value = 217
$ws.Cells.Item(72, 2).Value = 8606
$ws.Cells.Item(72, 3).Value = 6
$ws.Cells.Item(72, 4).Value = 8294
$ws.Cells.Item(72, 5).Value = 191
$ws.Cells.Item(76, 2).Value = 7191
$ws.Cells.Item(76, 3).Value = 19
$ws.Cells.Item(76, 5).Value = 264
$ws.Cells.Item(77, 2).Value = 6552
$ws.Cells.Item(77, 3).Value = 141
$ws.Cells.Item(77, 4).Value = 900
$ws.Cells.Item(77, 5).Value = 5503
$ws.Cells.Item(77, 7).Value = 7
$ws.Cells.Item(77, 8).Value = 149
$ws.Cells.Item(78, 2).Value = 6354
$ws.Cells.Item(78, 3).Value = 121
$ws.Cells.Item(78, 4).Value = 4193
$ws.Cells.Item(78, 5).Value = 2063
$ws.Cells.Item(78, 7).Value = 4
$ws.Cells.Item(78, 8).Value = 98
$ws.Cells.Item(82, 1).Value = "El Salvador"
$ws.Cells.Item(82, 2).Value = 5517
$ws.Cells.Item(82, 3).Value = 181
$ws.Cells.Item(82, 4).Value = 3291
$ws.Cells.Item(82, 5).Value = 2093
$ws.Cells.Item(82, 7).Value = 7
$ws.Cells.Item(82, 8).Value = 133
$ws.Cells.Item(83, 1).Value = "Kenia"
$ws.Cells.Item(83, 2).Value = 5384
$ws.Cells.Item(83, 4).Value = 1857
$ws.Cells.Item(83, 5).Value = 3395
$ws.Cells.Item(83, 7).Value = 0
$ws.Cells.Item(83, 8).Value = 132
$ws.Cells.Item(105, 2).Value = 2269
$ws.Cells.Item(105, 3).Value = 77
$ws.Cells.Item(105, 4).Value = 1298
$ws.Cells.Item(105, 5).Value = 920
$ws.Cells.Item(105, 7).Value = 2
$ws.Cells.Item(105, 8).Value = 51
$ws.Cells.Item(118, 1).Value = "Eslovenia"
$ws.Cells.Item(118, 2).Value = 1558
$ws.Cells.Item(118, 3).Value = 11
$ws.Cells.Item(118, 4).Value = 1376
$ws.Cells.Item(118, 5).Value = 73
$ws.Cells.Item(118, 8).Value = 109
$ws.Cells.Item(119, 1).Value = "Guinea-Bisau"
$ws.Cells.Item(119, 2).Value = 1556
$ws.Cells.Item(119, 4).Value = 191
$ws.Cells.Item(119, 5).Value = 1346
$ws.Cells.Item(119, 8).Value = 19
$ws.Cells.Item(121, 1).Value = "Estado de Palestina"
$ws.Cells.Item(121, 2).Value = 1514
$ws.Cells.Item(121, 3).Value = 132
$ws.Cells.Item(121, 4).Value = 446
$ws.Cells.Item(121, 5).Value = 1065
$ws.Cells.Item(121, 8).Value = 3
$ws.Cells.Item(122, 1).Value = "Zambia"
$ws.Cells.Item(122, 2).Value = 1497
$ws.Cells.Item(122, 4).Value = 1223
$ws.Cells.Item(122, 5).Value = 256
$ws.Cells.Item(122, 8).Value = 18
$ws.Cells.Item(124, 2).Value = 1197
$ws.Cells.Item(124, 3).Value = 3
$ws.Cells.Item(124, 4).Value = 1091
$ws.Cells.Item(130, 2).Value = 1059
$ws.Cells.Item(130, 3).Value = 3
$ws.Cells.Item(130, 4).Value = 919
$ws.Cells.Item(130, 5).Value = 73
$ws.Cells.Item(135, 2).Value = 941
$ws.Cells.Item(135, 3).Value = 7
$ws.Cells.Item(135, 4).Value = 830
$ws.Cells.Item(135, 5).Value = 58
$ws.Cells.Item(141, 2).Value = 833
$ws.Cells.Item(141, 3).Value = 12
$ws.Cells.Item(141, 4).Value = 761
$ws.Cells.Item(141, 5).Value = 72
$ws.Cells.Item(150, 2).Value = 670
$ws.Cells.Item(150, 3).Value = 2
$ws.Cells.Item(150, 4).Value = 632
$ws.Cells.Item(150, 5).Value = 29
$ws.Cells.Item(158, 2).Value = 353
$ws.Cells.Item(158, 3).Value = 1
$ws.Cells.Item(158, 4).Value = 330
$ws.Cells.Item(178, 1).Value = "Namibia"
$ws.Cells.Item(178, 2).Value = 105
$ws.Cells.Item(178, 3).Value = 3
$ws.Cells.Item(178, 4).Value = 21
$ws.Cells.Item(178, 5).Value = 84
$ws.Cells.Item(178, 8).Value = 0
$ws.Cells.Item(179, 1).Value = "Bahamas"
$ws.Cells.Item(179, 2).Value = 104
$ws.Cells.Item(179, 4).Value = 84
$ws.Cells.Item(179, 5).Value = 9
$ws.Cells.Item(179, 8).Value = 11
$ws.Cells.Item(186, 4).Value = 38
$ws.Cells.Item(186, 5).Value = 32
$ws.Cells.Item(202, 1).Value = "Fiyi"
$ws.Cells.Item(203, 1).Value = "Dominica"
$ws.Cells.Item(208, 1).Value = "Groenlandia"
$ws.Cells.Item(209, 1).Value = "Islas Malvinas"
$ws.Cells.Item(211, 1).Value = "Seychelles"
$ws.Cells.Item(211, 4).Value = 11
$ws.Cells.Item(211, 8).Value = 0
$ws.Cells.Item(212, 1).Value = "Montserrat"
$ws.Cells.Item(212, 4).Value = 10
$ws.Cells.Item(212, 8).Value = 1
